# Atualização de bases das ligas: swap a set of duplicate-looking fixture rows
# whose (HomeTeam, AwayTeam, odds, ...) data had been written into the wrong
# one of a same-date pair of rows. For each pair below, the full row content
# in columns B:AD (id/date stay in A and are shared/identical in C:D) is
# exchanged between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(23, 24),
    @(35, 36),
    @(38, 39),
    @(44, 45),
    @(65, 66),
    @(128, 129),
    @(167, 168),
    @(172, 173),
    @(187, 188),
    @(233, 234),
    @(291, 293)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
